$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so that
# numeric-looking strings (e.g. "311.41") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.443.63"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.512.51"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "311.41"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "98.68"
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("D7").Value = "0.562"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").Value = "35.06"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").Value = "0.0798"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "7.18"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").Value = "2.899.03"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "15.36"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").Value = "2.500.72"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "0.805"
$ws.Range("E17").Value = "  -3.92%  "
$ws.Range("D18").Value = "42.422.79"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "6.58"
$ws.Range("E19").Value = "  -4.27%  "
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").Value = "11.98"
$ws.Range("E21").Value = "  -3.66%  "
$ws.Range("D22").Value = "68.48"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "240.46"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").Value = "2.84"
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").Value = "1.98"
$ws.Range("E25").Value = "  -4.62%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "25.26"
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "37.97"
$ws.Range("E30").Value = "  -7.63%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "5.82"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "156.71"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").Value = "2.78"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "0.0782"
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").Value = "1.95"
$ws.Range("E37").Value = "  -6.18%  "
$ws.Range("D38").Value = "17.24"
$ws.Range("E38").Value = "  -6.80%  "
$ws.Range("D39").Value = "0.107"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "4.18"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").Value = "21.24"
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.007.46"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.23"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("D46").Value = "0.0294"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").Value = "9.14"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").Value = "2.751.02"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.187"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "78.41"
$ws.Range("E50").Value = "  -3.85%  "
$ws.Range("D51").Value = "99.84"
$ws.Range("E51").Value = "  -2.25%  "

# Restore the original (default) cell style now that the text values are set,
# so the saved workbook does not carry a spurious style/numberformat change.
$ws.Range("D2:E51").Style = "Normal"

